$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 95, shifting rows 95:123 down to 96:124
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new weekly record
$ws.Cells.Item(95, 1).Value = 8
$ws.Cells.Item(95, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(95, 3).Value = "Coquimbo"
$ws.Cells.Item(95, 4).Value = 44988
$ws.Cells.Item(95, 5).Value = 4
$ws.Cells.Item(95, 6).Value = 100112030
$ws.Cells.Item(95, 7).Value = "Poroto granado"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 440
$ws.Cells.Item(95, 11).Value = 37000
$ws.Cells.Item(95, 12).Value = 38000
$ws.Cells.Item(95, 13).Value = 37500
$ws.Cells.Item(95, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(95, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(95, 16).Value = 1500
$ws.Cells.Item(95, 17).Value = 25
$ws.Cells.Item(95, 18).Value = "Hortaliza"
